# Update QuantitativeEvaluation test-case results for UC3.4.4_TC1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Compilation success -> "no", with a note explaining why
$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = "Called wrong method"

# Runtime without error -> value cleared (no longer applicable)
$ws.Range("B6").ClearContents()

# Assertion validity -> value + note cleared
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Updated Code BLEU score + detail note
$ws.Range("B12").Value = 0.2630476915997347
$ws.Range("C12").Value = "{'codebleu': 0.2630476915997347, 'ngram_match_score': 0.10069921606275717, 'weighted_ngram_match_score': 0.11038491476272382, 'syntax_match_score': 0.5377358490566038, 'dataflow_match_score': 0.30337078651685395}"

# Reflect the new active selection
$ws.Range("B6").Select() | Out-Null
